$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 259
$ws.Range("F3").Value = 2782
$ws.Range("F7").Value = 2971
$ws.Range("F8").Value = 1896
$ws.Range("F10").Value = 70
$ws.Range("F11").Value = 2547
$ws.Range("F12").Value = 575
$ws.Range("F13").Value = 264
$ws.Range("F14").Value = 6
$ws.Range("F18").Value = 9495
$ws.Range("F19").Value = 63
$ws.Range("F21").Value = 7482
$ws.Range("F22").Value = 12012
$ws.Range("F26").Value = 377
$ws.Range("F27").Value = 579
$ws.Range("F28").Value = 2700
$ws.Range("F29").Value = 242
$ws.Range("F30").Value = 215
$ws.Range("F31").Value = 2692
$ws.Range("F32").Value = 999
$ws.Range("F34").Value = 60
$ws.Range("F35").Value = 55
$ws.Range("F36").Value = 4554
$ws.Range("F37").Value = 1088
$ws.Range("F38").Value = 32
$ws.Range("F41").Value = 568

$ws = $wb.Worksheets.Item(2)
$ws.Range("F3").Value = 72
$ws.Range("F21").Value = 17

$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 635
$ws.Range("F4").Value = 192

$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 635
$ws.Range("F3").Value = 259
$ws.Range("F5").Value = 2782
$ws.Range("F6").Value = 72
$ws.Range("F11").Value = 2971
$ws.Range("F13").Value = 1896
$ws.Range("F15").Value = 2547
$ws.Range("F17").Value = 575
$ws.Range("F18").Value = 264
$ws.Range("F19").Value = 6
$ws.Range("F22").Value = 9495
$ws.Range("F23").Value = 63
$ws.Range("F25").Value = 7482
$ws.Range("F26").Value = 12013
$ws.Range("F30").Value = 377
$ws.Range("F32").Value = 579
$ws.Range("F34").Value = 2700
$ws.Range("F36").Value = 242
$ws.Range("F37").Value = 215
$ws.Range("F38").Value = 60
$ws.Range("F39").Value = 55
$ws.Range("F40").Value = 4554
$ws.Range("F44").Value = 17
$ws.Range("F45").Value = 568
